$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values (e.g. "1.004") so they are not
# re-interpreted as numbers by Excel's input parser.
$numericLookingCells = @("D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D20","D22","D24","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values scraped on Sat Apr 29 14:41:31 UTC 2023
$ws.Range("D2").Value = '29.489.06'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '1.923.33'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("E4").Value = '  +0.44%  '
$ws.Range("D5").Value = '325.49'
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("D7").Value = '0.4841'
$ws.Range("E7").Value = '  +3.14%  '
$ws.Range("D8").Value = '0.4100'
$ws.Range("E8").Value = '  +2.07%  '
$ws.Range("D9").Value = '0.08182'
$ws.Range("E9").Value = '  +2.37%  '
$ws.Range("D10").Value = '1.026'
$ws.Range("E10").Value = '  +3.44%  '
$ws.Range("D11").Value = '23.88'
$ws.Range("E11").Value = '  +6.66%  '
$ws.Range("D12").Value = '1.919.51'
$ws.Range("E12").Value = '  +2.43%  '
$ws.Range("D13").Value = '6.052'
$ws.Range("E13").Value = '  +3.57%  '
$ws.Range("D14").Value = '7.224'
$ws.Range("E14").Value = '  +2.84%  '
$ws.Range("D15").Value = '91.57'
$ws.Range("E15").Value = '  +3.12%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '0.06763'
$ws.Range("E17").Value = '  +2.49%  '
$ws.Range("E19").Value = '  +2.15%  '
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = '29.514.84'
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("D22").Value = '5.635'
$ws.Range("E22").Value = '  +2.64%  '
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("D24").Value = '2.179'
$ws.Range("E24").Value = '  -0.78%  '
$ws.Range("D25").Value = '2.142.36'
$ws.Range("E25").Value = '  +2.02%  '
$ws.Range("D26").Value = '6.697'
$ws.Range("E26").Value = '  +9.71%  '
$ws.Range("D27").Value = '156.79'
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("D28").Value = '20.12'
$ws.Range("E28").Value = '  +2.54%  '
$ws.Range("D29").Value = '2.129'
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("E30").Value = '  +3.32%  '
$ws.Range("D31").Value = '1.030'
$ws.Range("E31").Value = '  -1.74%  '
$ws.Range("D32").Value = '0.09578'
$ws.Range("E32").Value = '  +1.71%  '
$ws.Range("D33").Value = '5.548'
$ws.Range("E33").Value = '  +4.04%  '
$ws.Range("D34").Value = '3.566'
$ws.Range("E34").Value = '  +0.71%  '
$ws.Range("D35").Value = '1.395'
$ws.Range("E35").Value = '  +0.29%  '
$ws.Range("D36").Value = '0.02291'
$ws.Range("E36").Value = '  +2.60%  '
$ws.Range("D37").Value = '0.06154'
$ws.Range("E37").Value = '  +1.54%  '
$ws.Range("D38").Value = '1.182'
$ws.Range("E38").Value = '  +0.76%  '
$ws.Range("D39").Value = '0.5993'
$ws.Range("E39").Value = '  +3.36%  '
$ws.Range("D40").Value = '10.82'
$ws.Range("E40").Value = '  +8.20%  '
$ws.Range("D41").Value = '8.019'
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("E42").Value = '  +2.54%  '
$ws.Range("D43").Value = '2.458'
$ws.Range("E43").Value = '  -0.42%  '
$ws.Range("D44").Value = '1.280'
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").Value = '12.47'
$ws.Range("E46").Value = '  +3.15%  '
$ws.Range("D47").Value = '0.5599'
$ws.Range("E47").Value = '  +2.54%  '
$ws.Range("D48").Value = '1.966'
$ws.Range("E48").Value = '  +3.71%  '
$ws.Range("D49").Value = '116.97'
$ws.Range("E49").Value = '  +3.19%  '
$ws.Range("D50").Value = '2.434'
$ws.Range("E50").Value = '  +4.16%  '
$ws.Range("D51").Value = '72.96'
$ws.Range("E51").Value = '  +3.11%  '

# Restore default (Normal) style on those cells so only the value text changed,
# matching the source workbook which keeps them unstyled text cells.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
